{"js": "// Replace each exact cell/date string with its updated value.\n// Each source string in this worksheet-style table is unique, so a single\n// search+replace per pair is sufficient and keeps run formatting (rFonts/sz) intact.\nconst body = context.document.body;\n\nconst replacements = [\n  [\"2025-12-13 Saturday\", \"2025-12-14 Sunday\"],\n  [\"141\u00d76=846\", \"444\u00d76=2664\"],\n  [\"578\u00d73=1734\", \"656\u00d78=5248\"],\n  [\"471\u00d78=3768\", \"116\u00d79=1044\"],\n  [\"419\u00d77=2933\", \"967\u00d79=8703\"],\n  [\"840\u00d79=7560\", \"678\u00d76=4068\"],\n  [\"603\u00d77=4221\", \"279\u00d76=1674\"],\n  [\"423\u00d75=2115\", \"535\u00d73=1605\"],\n  [\"476\u00d78=3808\", \"833\u00d77=5831\"],\n  [\"125\u00d72=250\", \"542\u00d73=1626\"],\n  [\"533\u00d77=3731\", \"186\u00d73=558\"],\n  [\"282\u00d76=1692\", \"134\u00d74=536\"],\n  [\"764\u00d75=3820\", \"749\u00d76=4494\"],\n  [\"787\u00d77=5509\", \"895\u00d75=4475\"],\n  [\"552\u00d77=3864\", \"785\u00d76=4710\"],\n  [\"345\u00d78=2760\", \"696\u00d77=4872\"],\n  [\"454\u00d78=3632\", \"368\u00d74=1472\"],\n  [\"298\u00d79=2682\", \"942\u00d76=5652\"],\n  [\"382\u00d73=1146\", \"880\u00d78=7040\"],\n  [\"406\u00d79=3654\", \"251\u00d76=1506\"],\n  [\"304\u00d73=912\", \"911\u00d72=1822\"],\n  [\"858\u00d72=1716\", \"793\u00d75=3965\"],\n  [\"588\u00d78=4704\", \"690\u00d74=2760\"],\n  [\"395\u00d75=1975\", \"114\u00d76=684\"],\n  [\"770\u00d79=6930\", \"501\u00d74=2004\"],\n  [\"366\u00d79=3294\", \"782\u00d78=6256\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Exact date/problem strings are unique in the document, so a plain\n# Find/Replace (wdReplaceAll) per pair swaps only the <w:t> text and\n# leaves every runs font/size formatting untouched.\n$replacements = @(\n    ,@(\"2025-12-13 Saturday\", \"2025-12-14 Sunday\")\n    ,@(\"141\u00d76=846\", \"444\u00d76=2664\")\n    ,@(\"578\u00d73=1734\", \"656\u00d78=5248\")\n    ,@(\"471\u00d78=3768\", \"116\u00d79=1044\")\n    ,@(\"419\u00d77=2933\", \"967\u00d79=8703\")\n    ,@(\"840\u00d79=7560\", \"678\u00d76=4068\")\n    ,@(\"603\u00d77=4221\", \"279\u00d76=1674\")\n    ,@(\"423\u00d75=2115\", \"535\u00d73=1605\")\n    ,@(\"476\u00d78=3808\", \"833\u00d77=5831\")\n    ,@(\"125\u00d72=250\", \"542\u00d73=1626\")\n    ,@(\"533\u00d77=3731\", \"186\u00d73=558\")\n    ,@(\"282\u00d76=1692\", \"134\u00d74=536\")\n    ,@(\"764\u00d75=3820\", \"749\u00d76=4494\")\n    ,@(\"787\u00d77=5509\", \"895\u00d75=4475\")\n    ,@(\"552\u00d77=3864\", \"785\u00d76=4710\")\n    ,@(\"345\u00d78=2760\", \"696\u00d77=4872\")\n    ,@(\"454\u00d78=3632\", \"368\u00d74=1472\")\n    ,@(\"298\u00d79=2682\", \"942\u00d76=5652\")\n    ,@(\"382\u00d73=1146\", \"880\u00d78=7040\")\n    ,@(\"406\u00d79=3654\", \"251\u00d76=1506\")\n    ,@(\"304\u00d73=912\", \"911\u00d72=1822\")\n    ,@(\"858\u00d72=1716\", \"793\u00d75=3965\")\n    ,@(\"588\u00d78=4704\", \"690\u00d74=2760\")\n    ,@(\"395\u00d75=1975\", \"114\u00d76=684\")\n    ,@(\"770\u00d79=6930\", \"501\u00d74=2004\")\n    ,@(\"366\u00d79=3294\", \"782\u00d78=6256\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
